$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.332.52"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "1.923.74"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8109"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3263"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.22"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07248"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7948"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08119"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.49%  "
$ws.Range("D13").Value = "1.913.34"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.450"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.83%  "
$ws.Range("E15").Value = "  +1.46%  "
$ws.Range("D16").Value = "30.330.96"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.101"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "250.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007872"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.36%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.280"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +20.87%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.175.81"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1653"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +18.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.544"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.167"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.388"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.556"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.357"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05767"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.142"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.20%  "
$ws.Range("E35").Value = "  +3.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7495"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.18%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.750"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("B38").Value = "Frax"
$ws.Range("C38").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9985"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.819"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4518"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8586"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.936"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.89%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.039.87"
$ws.Range("E47").Value = "  +5.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.144"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.675"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.946"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.54%  "
